# Account.xlsx - "Account AutoIncrement 옵션 누락" (missing AutoIncrement option)
#
# The "Id" row's Key column ("pk") was missing a note that the primary key
# is autogenerated (auto-increment). Update it to "pk, autogenerated".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the "Id" row; column F is "Key".
$ws.Range("F2").Value = "pk, autogenerated"

# Leave the selection where the author ended up after making the edit.
$ws.Range("F3").Select()
